$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor-expressing cell counts and resulting recalculated
# expression metrics (per Dr Hou's advice on natmi analysis).
$data = @{
    2 = @(3, 1, 15.84401633333333, 47.532049, 0.7212145038223592, 0.7256686596529557, 3, 1, 56.01527666666667, 168.04583, 0.4883616647765734, 0.5237161609491596, 887.5069584228522, 7987.56262580567, 0.3522135157476977, 0.3800444045545683)
    3 = @(3, 1, 15.84401633333333, 47.532049, 0.7212145038223592, 0.7256686596529557, 3, 1, 9.149484, 27.448452, 0.0797685471532371, 0.08554331818550501, 144.964573937572, 1304.681165438148, 0.05753023315575236, 0.06207610504994173)
    4 = @(3, 1, 15.84401633333333, 47.532049, 0.7212145038223592, 0.7256686596529557, 3, 1, 12.129045, 36.387135, 0.1057454494708373, 0.1134007945935868, 192.172787087735, 1729.555083789615, 0.07626515187158231, 0.08229140261630825)
    5 = @(3, 1, 15.84401633333333, 47.532049, 0.7212145038223592, 0.7256686596529557, 3, 1, 14.177359, 42.532077, 0.1236033999184949, 0.1325515550349214, 224.6263075595303, 2021.636768035773, 0.08914456474297394, 0.09618850927710643)
    6 = @(3, 1, 15.84401633333333, 47.532049, 0.7212145038223592, 0.7256686596529557, 2, 1, 23.229232, 46.458464, 0.2025209386808572, 0.1447881712368271, 368.0443312187894, 2208.265987312736, 0.1460610383043529, 0.105068238155031)
    7 = @(3, 1, 2.951312666666666, 8.853938, 0.1343427989301267, 0.1351724879588164, 3, 1, 56.01527666666667, 168.04583, 0.4883616647765734, 0.5237161609491596, 165.3185955531711, 1487.86735997854, 0.06560787293626115, 0.07079201645973783)
    8 = @(3, 1, 2.951312666666666, 8.853938, 0.1343427989301267, 0.1351724879588164, 3, 1, 9.149484, 27.448452, 0.0797685471532371, 0.08554331818550501, 27.002988022664, 243.026892203976, 0.01071632989115566, 0.01156310314738737)
    9 = @(3, 1, 2.951312666666666, 8.853938, 0.1343427989301267, 0.1351724879588164, 3, 1, 12.129045, 36.387135, 0.1057454494708373, 0.1134007945935868, 35.79660414307, 322.16943728763, 0.01420613965603658, 0.01532866754172182)
    10 = @(3, 1, 2.951312666666666, 8.853938, 0.1343427989301267, 0.1351724879588164, 3, 1, 14.177359, 42.532077, 0.1236033999184949, 0.1325515550349214, 41.84181919658067, 376.576372769226, 0.0166052267023304, 0.01791732347688031)
    11 = @(3, 1, 2.951312666666666, 8.853938, 0.1343427989301267, 0.1351724879588164, 2, 1, 23.229232, 46.458464, 0.2025209386808572, 0.1447881712368271, 68.55672663853866, 411.340359831232, 0.02720722974434292, 0.01957137733308906)
    12 = @(3, 1, 2.768664, 8.305992, 0.1260286906426543, 0.1268070324872419, 3, 1, 56.01527666666667, 168.04583, 0.4883616647765734, 0.5237161609491596, 155.08747995704, 1395.78731961336, 0.0615475811718584, 0.06641089223557367)
    13 = @(3, 1, 2.768664, 8.305992, 0.1260286906426543, 0.1268070324872419, 3, 1, 9.149484, 27.448452, 0.0797685471532371, 0.08554331818550501, 25.331846969376, 227.986622724384, 0.0100531255521893, 0.0108474943282158)
    14 = @(3, 1, 2.768664, 8.305992, 0.1260286906426543, 0.1268070324872419, 3, 1, 12.129045, 36.387135, 0.1057454494708373, 0.1134007945935868, 33.58125024588, 302.23125221292, 0.01332696053822859, 0.014380018244108)
    15 = @(3, 1, 2.768664, 8.305992, 0.1260286906426543, 0.1268070324872419, 3, 1, 14.177359, 42.532077, 0.1236033999184949, 0.1325515550349214, 39.252343478376, 353.271091305384, 0.01557757465070827, 0.01680846934554771)
    16 = @(3, 1, 2.768664, 8.305992, 0.1260286906426543, 0.1268070324872419, 2, 1, 23.229232, 46.458464, 0.2025209386808572, 0.1447881712368271, 64.313938386048, 385.883630316288, 0.02552344872966971, 0.01836015833379668)
    17 = @(2, 1, 0.4045285, 0.809057, 0.01841400660485959, 0.01235181990098599, 3, 1, 56.01527666666667, 168.04583, 0.4883616647765734, 0.5237161609491596, 22.65977584705167, 135.95865508231, 0.008992694920756049, 0.006468847699279814)
    18 = @(2, 1, 0.4045285, 0.809057, 0.01841400660485959, 0.01235181990098599, 3, 1, 9.149484, 27.448452, 0.0797685471532371, 0.08554331818550501, 3.701227038294, 22.207362229764, 0.001468858554139762, 0.001056615659960098)
    19 = @(2, 1, 0.4045285, 0.809057, 0.01841400660485959, 0.01235181990098599, 3, 1, 12.129045, 36.387135, 0.1057454494708373, 0.1134007945935868, 4.9065443802825, 29.439266281695, 0.001947197404989845, 0.00140070619144869)
    20 = @(2, 1, 0.4045285, 0.809057, 0.01841400660485959, 0.01235181990098599, 3, 1, 14.177359, 42.532077, 0.1236033999184949, 0.1325515550349214, 5.735145770231501, 34.410874621389, 0.002276033822482267, 0.001637252935386983)
    21 = @(2, 1, 0.4045285, 0.809057, 0.01841400660485959, 0.01235181990098599, 2, 1, 23.229232, 46.458464, 0.2025209386808572, 0.1447881712368271, 9.396886377112, 37.587545508448, 0.003729221902491669, 0.001788397414910409)
}

$columns = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
